$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 2.27
$ws.Range("H2").Value = 2.75
$ws.Range("I2").Value = 3.6
$ws.Range("J2").Value = 2.95
$ws.Range("K2").Value = 1.85
$ws.Range("L2").Value = 4.3
$ws.Range("N2").Value = 5.2
$ws.Range("O2").Value = 1.53
$ws.Range("P2").Value = 2.35
$ws.Range("Q2").Value = 2.52
$ws.Range("R2").Value = 1.47
$ws.Range("S2").Value = 4.5
$ws.Range("T2").Value = 1.16
$ws.Range("W2").Value = 2.05
$ws.Range("X2").Value = 1.7
$ws.Range("Y2").Value = 5.8
$ws.Range("Z2").Value = 9.75
$ws.Range("AA2").Value = 9.25
$ws.Range("AB2").Value = 24
$ws.Range("AC2").Value = 22
$ws.Range("AD2").Value = 40
$ws.Range("AE2").Value = 5.2
$ws.Range("AF2").Value = 5.5
$ws.Range("AG2").Value = 16.5
$ws.Range("AH2").Value = 110
$ws.Range("AJ2").Value = 7.8
$ws.Range("AK2").Value = 18
$ws.Range("AL2").Value = 13
$ws.Range("AM2").Value = 60
$ws.Range("G4").Value = 1.5
$ws.Range("H4").Value = 4.33
$ws.Range("I4").Value = 5.75
$ws.Range("J4").Value = 2
$ws.Range("K4").Value = 2.38
$ws.Range("L4").Value = 5.5
$ws.Range("Q4").Value = 1.65
$ws.Range("R4").Value = 2.2
$ws.Range("W4").Value = 1.8
$ws.Range("X4").Value = 1.91
$ws.Range("Z4").Value = 8
$ws.Range("AB4").Value = 11
$ws.Range("AC4").Value = 12
$ws.Range("AF4").Value = 8.5
$ws.Range("AG4").Value = 17
$ws.Range("AH4").Value = 51
$ws.Range("AJ4").Value = 17
$ws.Range("AK4").Value = 29
$ws.Range("AL4").Value = 17
$ws.Range("AN4").Value = 41
$ws.Range("S5").Value = 2
$ws.Range("T5").Value = 1.73
$ws.Range("AD5").Value = 21
$ws.Range("AE5").Value = 21
$ws.Range("J6").Value = 1.53
$ws.Range("L6").Value = 9
$ws.Range("AB6").Value = 7.5
$ws.Range("AH6").Value = 67
$ws.Range("AI6").Value = 600
